# Update Means sheet (sheet1): replace #NUM! errors in D2:G10 with computed values
$wb = $excel.ActiveWorkbook
$wsMeans = $wb.Worksheets.Item("Means")

$meansData = @{
    2  = @(99, 97, 96, 90)
    3  = @(0, 0.29, 0.96, 6.2)
    4  = @(0.89, 2.7, 2.9, 3.9)
    5  = @(0.45, 0.91, 0.9, 0.89)
    6  = @(48, 47, 44, 49)
    7  = @(10, 11, 11, 9)
    8  = @(5.5, 7.4, 5.9, 9.1)
    9  = @(30, 30, 30, 31)
    10 = @(0.4, 0.4, 0.4, 0.39)
}

foreach ($row in $meansData.Keys) {
    $vals = $meansData[$row]
    $wsMeans.Cells.Item($row, 4).Value = $vals[0]   # D
    $wsMeans.Cells.Item($row, 5).Value = $vals[1]   # E
    $wsMeans.Cells.Item($row, 6).Value = $vals[2]   # F
    $wsMeans.Cells.Item($row, 7).Value = $vals[3]   # G
}

# Update Standard Deviations sheet (sheet2): update E2:G10 (D column stays 0)
$wsSD = $wb.Worksheets.Item("Standard Deviations")

$sdData = @{
    2  = @(3.2, 5.1, 13)
    3  = @(0.67, 2.8, 11)
    4  = @(2.8, 2.8, 4.3)
    5  = @(1.1, 1, 1.7)
    6  = @(12, 9.4, 21)
    7  = @(9.1, 10, 8.8)
    8  = @(5.4, 4.5, 9.1)
    9  = @(0, 0, 2.8)
    10 = @(0.000000000000000028, 0.000000000000000019, 0.027)
}

foreach ($row in $sdData.Keys) {
    $vals = $sdData[$row]
    $wsSD.Cells.Item($row, 5).Value = $vals[0]   # E
    $wsSD.Cells.Item($row, 6).Value = $vals[1]   # F
    $wsSD.Cells.Item($row, 7).Value = $vals[2]   # G
}
